$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data (row 2), shifting existing data rows down
$ws.Range("A2:C5").Insert()
$ws.Range("A2:C5").ClearFormats()

# Fill in the newly inserted rows with data
$ws.Cells.Item(2, 1).Value = 2.563363254070282
$ws.Cells.Item(2, 2).Value = -5.564052700996399
$ws.Cells.Item(2, 3).Value = -4.925167679786682
$ws.Cells.Item(3, 1).Value = 2.571200489997864
$ws.Cells.Item(3, 2).Value = -5.45090651512146
$ws.Cells.Item(3, 3).Value = -4.94497549533844
$ws.Cells.Item(4, 1).Value = 2.582025349140167
$ws.Cells.Item(4, 2).Value = -5.429405391216278
$ws.Cells.Item(4, 3).Value = -4.891633093357086
$ws.Cells.Item(5, 1).Value = 2.521161556243896
$ws.Cells.Item(5, 2).Value = -5.436496257781982
$ws.Cells.Item(5, 3).Value = -4.74793529510498

# Append new rows of data at the bottom (rows 26-31)
$ws.Cells.Item(26, 1).Value = -2.611050009727472
$ws.Cells.Item(26, 2).Value = -2.573673054575919
$ws.Cells.Item(26, 3).Value = -8.294337868690496
$ws.Cells.Item(27, 1).Value = 0.61596310138703
$ws.Cells.Item(27, 2).Value = -2.870795279741297
$ws.Cells.Item(27, 3).Value = -8.901223957538587
$ws.Cells.Item(28, 1).Value = 2.079445004463198
$ws.Cells.Item(28, 2).Value = -5.368536770343783
$ws.Cells.Item(28, 3).Value = -3.524431616067886
$ws.Cells.Item(29, 1).Value = 0.3216586112976074
$ws.Cells.Item(29, 2).Value = -3.676267147064209
$ws.Cells.Item(29, 3).Value = -3.865855693817138
$ws.Cells.Item(30, 1).Value = 1.706132471561434
$ws.Cells.Item(30, 2).Value = -4.47040206193924
$ws.Cells.Item(30, 3).Value = -5.197765350341799
$ws.Cells.Item(31, 1).Value = 1.816254138946533
$ws.Cells.Item(31, 2).Value = -3.920204520225524
$ws.Cells.Item(31, 3).Value = -5.85852086544037
